# Weekly crime data refresh: volume/date header bump, updated statistics,
# and an extra footer row inserted before the "Prepared by" block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header text: Volume number and the reporting week date range
# ---------------------------------------------------------------------------
$ws.Range("A8").Value  = "Volume 31   Number  51"
$ws.Range("C9").Value  = "Report Covering the Week  12/16/2024  Through  12/22/2024"

# ---------------------------------------------------------------------------
# 2. Insert a blank row before the old row 56 so the "Prepared by" /
#    "NYPD CompStat Unit" block moves from rows 56-57 down to rows 57-58.
# ---------------------------------------------------------------------------
$ws.Rows.Item(56).Insert()

# ---------------------------------------------------------------------------
# 3. Updated crime-statistics figures (rows 16-31)
# ---------------------------------------------------------------------------
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -20
$ws.Range("F16").Value = 13
$ws.Range("H16").Value = -23.529411764705
$ws.Range("I16").Value = 228
$ws.Range("J16").Value = 255
$ws.Range("K16").Value = -10.588235294117
$ws.Range("L16").Value = 20
$ws.Range("M16").Value = 21.925133689839
$ws.Range("N16").Value = -77.2

$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -40
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = 17.647058823529
$ws.Range("I17").Value = 255
$ws.Range("J17").Value = 234
$ws.Range("K17").Value = 8.974358974358
$ws.Range("L17").Value = 18.055555555555
$ws.Range("M17").Value = 102.380952380952
$ws.Range("N17").Value = -16.938110749185

$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 17
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = -29.166666666666
$ws.Range("I18").Value = 238
$ws.Range("J18").Value = 225
$ws.Range("K18").Value = 5.777777777777
$ws.Range("L18").Value = 37.572254335260
$ws.Range("M18").Value = -5.555555555555
$ws.Range("N18").Value = -84.403669724770

$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 22
$ws.Range("E19").Value = -63.636363636363
$ws.Range("F19").Value = 45
$ws.Range("G19").Value = 66
$ws.Range("H19").Value = -31.818181818181
$ws.Range("I19").Value = 701
$ws.Range("J19").Value = 736
$ws.Range("K19").Value = -4.755434782608
$ws.Range("L19").Value = 3.851851851851
$ws.Range("M19").Value = 55.777777777777
$ws.Range("N19").Value = -22.024471635150

$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = -66.666666666666
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 22
$ws.Range("H20").Value = -40.909090909090
$ws.Range("I20").Value = 249
$ws.Range("J20").Value = 304
$ws.Range("K20").Value = -18.092105263157
$ws.Range("L20").Value = 6.866952789699
$ws.Range("M20").Value = 29.6875
$ws.Range("N20").Value = -87.982625482625

$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 43
$ws.Range("E21").Value = -53.488372093023
$ws.Range("F21").Value = 108
$ws.Range("G21").Value = 146
$ws.Range("H21").Value = -26.027397260274
$ws.Range("I21").Value = 1693
$ws.Range("J21").Value = 1778
$ws.Range("K21").Value = -4.780652418447
$ws.Range("L21").Value = 12.267904509283
$ws.Range("M21").Value = 38.091353996737
$ws.Range("N21").Value = -70.940611053896

$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -50
$ws.Range("F22").Value = 8
$ws.Range("H22").Value = 33.333333333333
$ws.Range("I22").Value = 62
$ws.Range("J22").Value = 85
$ws.Range("K22").Value = -27.058823529411
$ws.Range("L22").Value = -20.512820512820
$ws.Range("M22").Value = 26.530612244898

$ws.Range("C24").Value = 39
$ws.Range("D24").Value = 42
$ws.Range("E24").Value = -7.142857142857
$ws.Range("F24").Value = 175
$ws.Range("H24").Value = 1.156069364161
$ws.Range("I24").Value = 2026
$ws.Range("J24").Value = 2036
$ws.Range("K24").Value = -0.491159135559
$ws.Range("L24").Value = 32.852459016393
$ws.Range("M24").Value = 119.501625135428

$ws.Range("C25").Value = 19
$ws.Range("D25").Value = 31
$ws.Range("E25").Value = -38.709677419354
$ws.Range("F25").Value = 112
$ws.Range("G25").Value = 109
$ws.Range("H25").Value = 2.752293577981
$ws.Range("I25").Value = 1366
$ws.Range("J25").Value = 1264
$ws.Range("K25").Value = 8.069620253164
$ws.Range("L25").Value = 122.838499184339

$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 12
$ws.Range("E26").Value = -25
$ws.Range("F26").Value = 43
$ws.Range("H26").Value = -2.272727272727
$ws.Range("I26").Value = 550
$ws.Range("J26").Value = 515
$ws.Range("K26").Value = 6.796116504854
$ws.Range("L26").Value = 1.851851851851
$ws.Range("M26").Value = 11.336032388664

# Row 27: D27 and E27 switch from numbers to the "no data" text markers
# ("0" / "***.*") already used elsewhere in the sheet - copy an existing
# cell with that exact value+style so the shared-string type matches.
$ws.Range("C27").Copy($ws.Range("D27"))
$ws.Range("E14").Copy($ws.Range("E27"))

$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 4
$ws.Range("H28").Value = -50
$ws.Range("I28").Value = 64
$ws.Range("J28").Value = 85
$ws.Range("K28").Value = -24.705882352941
$ws.Range("L28").Value = -31.914893617021

# Row 31: F31 switches from the "0" text marker to an actual number.
$ws.Range("I31").Copy($ws.Range("F31"))
$ws.Range("F31").Value = 1
$ws.Range("I31").Value = 10
$ws.Range("K31").Value = 11.111111111111
$ws.Range("L31").Value = 11.111111111111

# Row 33: F33, G33 and H33 switch from numbers to the "no data" text
# markers ("0" / "***.*"), matching the same pattern used in row 27.
$ws.Range("C33").Copy($ws.Range("F33"))
$ws.Range("C33").Copy($ws.Range("G33"))
$ws.Range("E33").Copy($ws.Range("H33"))

$excel.CutCopyMode = $false
